# Apply updated crypto price / volume values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.370.96'
$ws.Range("E2").Value = '  +2.91%  '
$ws.Range("D3").Value = '2.306.54'
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Formula = "'" + '310.92'
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").Formula = "'" + '103.14'
$ws.Range("E6").Value = '  +6.18%  '
$ws.Range("D7").Formula = "'" + '0.531'
$ws.Range("E7").Value = '  +1.62%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Formula = "'" + '0.529'
$ws.Range("E9").Value = '  +8.38%  '
$ws.Range("D10").Formula = "'" + '36.23'
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("D11").Formula = "'" + '0.0811'
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("D12").Formula = "'" + '51.91'
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").Formula = "'" + '7.04'
$ws.Range("E14").Value = '  +3.45%  '
$ws.Range("D15").Value = '2.665.79'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D17").Value = '2.314.00'
$ws.Range("E17").Value = '  +2.49%  '
$ws.Range("D18").Formula = "'" + '0.809'
$ws.Range("E18").Value = '  +2.62%  '
$ws.Range("D19").Value = '43.276.80'
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("D20").Formula = "'" + '12.16'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").Value = '0.0₃0927'
$ws.Range("E21").Value = '  +2.94%  '
$ws.Range("E22").Value = '  +3.62%  '
$ws.Range("D23").Formula = "'" + '68.10'
$ws.Range("D24").Formula = "'" + '241.78'
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("E25").Value = '  +2.60%  '
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Formula = "'" + '24.76'
$ws.Range("E28").Value = '  +5.37%  '
$ws.Range("E29").Value = '  +8.11%  '
$ws.Range("D30").Formula = "'" + '36.82'
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("D31").Formula = "'" + '9.67'
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("D32").Formula = "'" + '168.04'
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Formula = "'" + '18.04'
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("E36").Value = '  +6.05%  '
$ws.Range("D37").Formula = "'" + '0.0742'
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("D38").Formula = "'" + '3.04'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("E39").Value = '  +3.64%  '
$ws.Range("E40").Value = '  +2.64%  '
$ws.Range("D41").Formula = "'" + '4.46'
$ws.Range("E41").Value = '  +8.09%  '
$ws.Range("D43").Formula = "'" + '2.58'
$ws.Range("E43").Value = '  +12.58%  '
$ws.Range("D44").Formula = "'" + '0.0295'
$ws.Range("E44").Value = '  +5.56%  '
$ws.Range("D45").Value = '1.986.06'
$ws.Range("E45").Value = '  +1.78%  '
$ws.Range("D46").Formula = "'" + '19.01'
$ws.Range("E46").Value = '  +0.53%  '
$ws.Range("D47").Formula = "'" + '3.00'
$ws.Range("E47").Value = '  +2.87%  '
$ws.Range("D48").Formula = "'" + '9.96'
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("D49").Formula = "'" + '55.84'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  +9.23%  '
